$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.095.18'
$ws.Range('D2').Style = $style_D2

$style_E2 = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('E2').Style = $style_E2

$style_D3 = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.825.33'
$ws.Range('D3').Style = $style_D3

$style_D4 = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').Style = $style_D4

$style_E4 = $ws.Range('E4').Style
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('E4').Style = $style_E4

$style_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.51'
$ws.Range('D5').Style = $style_D5

$style_E5 = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('E5').Style = $style_E5

$style_D6 = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('D6').Style = $style_D6

$style_E6 = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('E6').Style = $style_E6

$style_D7 = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4684'
$ws.Range('D7').Style = $style_D7

$style_E7 = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E7').Style = $style_E7

$style_D8 = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3649'
$ws.Range('D8').Style = $style_D8

$style_E8 = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('E8').Style = $style_E8

$style_D9 = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07384'
$ws.Range('D9').Style = $style_D9

$style_E9 = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('E9').Style = $style_E9

$style_D10 = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8785'
$ws.Range('D10').Style = $style_D10

$style_E10 = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('E10').Style = $style_E10

$style_D11 = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.25'
$ws.Range('D11').Style = $style_D11

$style_E11 = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('E11').Style = $style_E11

$style_D12 = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.895.38'
$ws.Range('D12').Style = $style_D12

$style_E12 = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.50%  '
$ws.Range('E12').Style = $style_E12

$style_D13 = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07481'
$ws.Range('D13').Style = $style_D13

$style_E13 = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.54%  '
$ws.Range('E13').Style = $style_E13

$style_D14 = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.370'
$ws.Range('D14').Style = $style_D14

$style_E14 = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('E14').Style = $style_E14

$style_E15 = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('E15').Style = $style_E15

$style_D16 = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.524'
$ws.Range('D16').Style = $style_D16

$style_E16 = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E16').Style = $style_E16

$style_E17 = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('E17').Style = $style_E17

$style_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008715'
$ws.Range('D18').Style = $style_D18

$style_E18 = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('E18').Style = $style_E18

$style_E19 = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('E19').Style = $style_E19

$style_D20 = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.429.89'
$ws.Range('D20').Style = $style_D20

$style_E20 = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('E20').Style = $style_E20

$style_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.60'
$ws.Range('D21').Style = $style_D21

$style_E21 = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('E21').Style = $style_E21

$style_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.232'
$ws.Range('D22').Style = $style_D22

$style_E22 = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('E22').Style = $style_E22

$style_D23 = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('D23').Style = $style_D23

$style_E23 = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E23').Style = $style_E23

$style_D24 = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.076.30'
$ws.Range('D24').Style = $style_D24

$style_E24 = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('E24').Style = $style_E24

$style_E25 = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E25').Style = $style_E25

$style_D26 = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.40'
$ws.Range('D26').Style = $style_D26

$style_E26 = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('E26').Style = $style_E26

$style_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.50'
$ws.Range('D27').Style = $style_D27

$style_E27 = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('E27').Style = $style_E27

$style_D28 = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.135'
$ws.Range('D28').Style = $style_D28

$style_E28 = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('E28').Style = $style_E28

$style_D29 = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.161'
$ws.Range('D29').Style = $style_D29

$style_E29 = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('E29').Style = $style_E29

$style_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.29'
$ws.Range('D30').Style = $style_D30

$style_E30 = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('E30').Style = $style_E30

$style_D31 = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08919'
$ws.Range('D31').Style = $style_D31

$style_E31 = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E31').Style = $style_E31

$style_D32 = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7436'
$ws.Range('D32').Style = $style_D32

$style_E32 = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('E32').Style = $style_E32

$style_D33 = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.161'
$ws.Range('D33').Style = $style_D33

$style_E33 = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E33').Style = $style_E33

$style_D34 = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.506'
$ws.Range('D34').Style = $style_D34

$style_E34 = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('E34').Style = $style_E34

$style_D35 = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.942'
$ws.Range('D35').Style = $style_D35

$style_E35 = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E35').Style = $style_E35

$style_E36 = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E36').Style = $style_E36

$style_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.535'
$ws.Range('D37').Style = $style_D37

$style_E37 = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.53%  '
$ws.Range('E37').Style = $style_E37

$style_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.090'
$ws.Range('D38').Style = $style_D38

$style_E38 = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E38').Style = $style_E38

$style_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05289'
$ws.Range('D39').Style = $style_D39

$style_E39 = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('E39').Style = $style_E39

$style_B40 = $ws.Range('B40').Style
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('B40').Style = $style_B40

$style_C40 = $ws.Range('C40').Style
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C40').Style = $style_C40

$style_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.337'
$ws.Range('D40').Style = $style_D40

$style_E40 = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('E40').Style = $style_E40

$style_B41 = $ws.Range('B41').Style
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'VeChain'
$ws.Range('B41').Style = $style_B41

$style_C41 = $ws.Range('C41').Style
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C41').Style = $style_C41

$style_D41 = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01931'
$ws.Range('D41').Style = $style_D41

$style_E41 = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.80%  '
$ws.Range('E41').Style = $style_E41

$style_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.934'
$ws.Range('D42').Style = $style_D42

$style_E42 = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('E42').Style = $style_E42

$style_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5245'
$ws.Range('D43').Style = $style_D43

$style_E43 = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('E43').Style = $style_E43

$style_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1639'
$ws.Range('D44').Style = $style_D44

$style_E44 = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('E44').Style = $style_E44

$style_D45 = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.366'
$ws.Range('D45').Style = $style_D45

$style_E45 = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('E45').Style = $style_E45

$style_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4892'
$ws.Range('D46').Style = $style_D46

$style_E46 = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('E46').Style = $style_E46

$style_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.40'
$ws.Range('D47').Style = $style_D47

$style_E47 = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('E47').Style = $style_E47

$style_E48 = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('E48').Style = $style_E48

$style_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.27'
$ws.Range('D49').Style = $style_D49

$style_E49 = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('E49').Style = $style_E49

$style_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06262'
$ws.Range('D51').Style = $style_D51

$style_E51 = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.52%  '
$ws.Range('E51').Style = $style_E51
